# Bugfixes: descontar tipos, 1T+1A=expulsion, modal continuar, tracking faltas
#
# The sheet is a simple bug-tracking list:
#   column A / B rows 9-19  -> already fixed issues ("corregido" in B)
#   column A/B   rows 20-28 -> pending issues (numbered 1-9 in A, text in B)
#   column A/B   rows 30-33 -> more pending issues (numbered 1-4 in A, text in B)
#   column A     rows 34-36 -> more pending issues (numbered 15-17 in A only)
#
# This change marks 6 of the pending issues as fixed (moves their text to
# column A in rows 20-25 and writes "corregido" into column B), and
# renumbers/repacks the remaining pending issues into rows 30-36 (text
# only in column A, no counters, no "corregido").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Read the current ("before") text of every pending-issue cell ----
# NOTE: use Value() with explicit parens to invoke the getter; the bare
# ".Value" property accessor (no parens) does not evaluate here.
$B20 = $ws.Range("B20").Value()
$B21 = $ws.Range("B21").Value()
$B22 = $ws.Range("B22").Value()
$B23 = $ws.Range("B23").Value()
$B24 = $ws.Range("B24").Value()
$B25 = $ws.Range("B25").Value()
$B26 = $ws.Range("B26").Value()
$B27 = $ws.Range("B27").Value()
$B28 = $ws.Range("B28").Value()
$B30 = $ws.Range("B30").Value()
$B31 = $ws.Range("B31").Value()
$B32 = $ws.Range("B32").Value()
$B33 = $ws.Range("B33").Value()

# ---- Clear out the old pending-issue area (rows 20-36, columns A:B) ----
$ws.Range("A20:B36").Clear() | Out-Null

# ---- Newly fixed issues: rows 20-25, text moves to column A, "corregido" in B ----
$ws.Range("A20").Value = $B30
$ws.Range("B20").Value = "corregido"

$ws.Range("A21").Value = $B31
$ws.Range("B21").Value = "corregido"

$ws.Range("A22").Value = $B33
$ws.Range("B22").Value = "corregido"

$ws.Range("A23").Value = $B20
$ws.Range("B23").Value = "corregido"

$ws.Range("A24").Value = $B26
$ws.Range("B24").Value = "corregido"

$ws.Range("A25").Value = $B27
$ws.Range("B25").Value = "corregido"

# ---- Remaining pending issues: rows 30-36, text only in column A ----
$ws.Range("A30").Value = $B21
$ws.Range("A31").Value = $B22
$ws.Range("A32").Value = $B23
$ws.Range("A33").Value = $B24
$ws.Range("A34").Value = $B25
$ws.Range("A35").Value = $B28
$ws.Range("A36").Value = $B32

# ---- Update the active selection shown when the workbook was saved ----
$ws.Range("B19:B25").Select() | Out-Null
